$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (F column) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1755
$wsExhibit.Range("F5").Value = 778
$wsExhibit.Range("F6").Value = 117
$wsExhibit.Range("F7").Value = 200

# Sheet "全部类型" (All Types) - update matching rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1755
$wsAll.Range("F6").Value = 778
$wsAll.Range("F7").Value = 117
$wsAll.Range("F8").Value = 200
